$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ref")

# Insert three new rows before the existing "absolute" row (row 14),
# shifting the rows below down by three.
$ws.Rows("14:16").Insert()

# New test case: expand-8
$ws.Range("B14").Value = "expand-8"
$ws.Range("C14").Value = "#P4:RR"

# New test case: expand-9
$ws.Range("B15").Value = "expand-9"
$ws.Range("C15").Value = "#P30:Q31:RR"

# New test case: expand-10
$ws.Range("B16").Value = "expand-10"
$ws.Range("C16").Value = "#M7:Q7:RRRRD"

# Make the "ref" sheet the active tab, with C16 selected (matches the
# author leaving off after typing the last new test case).
$ws.Activate() | Out-Null
$ws.Range("C16").Select() | Out-Null
